# Weekly update: insert a new daily price record as row 349 (Terminal La
# Palmera de La Serena - Zanahoria), shifting all subsequent rows down by
# one (old row 349 -> new row 350, ..., old row 430 -> new row 431).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 349; this pushes rows
# 349..430 down to 350..431 and extends the sheet dimension to R431.
$ws.Rows.Item(349).Insert()

# Populate the newly inserted row 349 with the new record. Most fields
# mirror the (now-shifted) row below it, except the date and the
# price/volume figures that differ for this new observation.
$ws.Range("A349").Value = 8
$ws.Range("B349").Value = "Terminal La Palmera de La Serena"
$ws.Range("C349").Value = "Coquimbo"
$ws.Range("D349").Value = 44889
$ws.Range("E349").Value = 4
$ws.Range("F349").Value = 100114013
$ws.Range("G349").Value = "Zanahoria"
$ws.Range("H349").Value = "Sin especificar"
$ws.Range("I349").Value = "Primera"
$ws.Range("J349").Value = 400
$ws.Range("K349").Value = 7000
$ws.Range("L349").Value = 8000
$ws.Range("M349").Value = 7500
$ws.Range("N349").Value = "`$/saco 20 kilos"
$ws.Range("O349").Value = "Provincia del Elquí"
$ws.Range("P349").Value = 375
$ws.Range("Q349").Value = 20
$ws.Range("R349").Value = "Hortaliza"
